$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.398934602737427
$ws.Range("B1").Value = 2.273723840713501
$ws.Range("C1").Value = 2.686756610870361
$ws.Range("D1").Value = 3.096260547637939
$ws.Range("E1").Value = 2.363765001296997
